$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 577145.5
$ws.Range("I28").Value = 757153.0600000001
$ws.Range("J28").Value = 1121.2
$ws.Range("K28").Value = 757153.0600000001
$ws.Range("L28").Value = 1121.2
$ws.Range("M28").Value = -756668.0600000001
$ws.Range("N28").Value = -2091.2

$ws.Range("H43").Value = 743.93335
$ws.Range("I43").Value = 651
$ws.Range("J43").Value = 758.2308
$ws.Range("K43").Value = 651
$ws.Range("L43").Value = 758.2308
$ws.Range("M43").Value = -582
$ws.Range("N43").Value = -896.2308

$ws.Range("H137").Value = 166668110
$ws.Range("I137").Value = 333334530
$ws.Range("J137").Value = 1700
$ws.Range("K137").Value = 1000003590
$ws.Range("L137").Value = 5100
$ws.Range("M137").Value = -1000001040
$ws.Range("N137").Value = -10200

$ws.Range("H138").Value = 6339318.5
$ws.Range("I138").Value = 1469161.4
$ws.Range("J138").Value = 8067438.5
$ws.Range("K138").Value = 4407484.199999999
$ws.Range("L138").Value = 24202315.5
$ws.Range("M138").Value = -4402344.199999999
$ws.Range("N138").Value = -24212595.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4144.353
$ws.Range("I61").Value = 2944.5715
$ws.Range("J61").Value = 4984.2
$ws.Range("K61").Value = 2944.5715
$ws.Range("L61").Value = 4984.2
$ws.Range("M61").Value = -2732.5715
$ws.Range("N61").Value = -5408.2

$ws.Range("H74").Value = 9932.4
$ws.Range("I74").Value = 2070.4285
$ws.Range("K74").Value = 2070.4285
$ws.Range("M74").Value = -1196.4285

$ws.Range("H77").Value = 9932.4
$ws.Range("I77").Value = 2070.4285
$ws.Range("K77").Value = 10352.1425
$ws.Range("M77").Value = -5984.1425

$ws.Range("H132").Value = 4003.7144
$ws.Range("I132").Value = 3464
$ws.Range("J132").Value = 4723.3335
$ws.Range("K132").Value = 10392
$ws.Range("L132").Value = 14170.0005
$ws.Range("M132").Value = -7862
$ws.Range("N132").Value = -19230.0005

$ws.Range("H136").Value = 4144.353
$ws.Range("I136").Value = 2944.5715
$ws.Range("J136").Value = 4984.2
$ws.Range("K136").Value = 8833.7145
$ws.Range("L136").Value = 14952.6
$ws.Range("M136").Value = -6283.7145
$ws.Range("N136").Value = -20052.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3728.9
$ws.Range("I134").Value = 2712.3635
$ws.Range("J134").Value = 4971.3335
$ws.Range("K134").Value = 8137.0905
$ws.Range("L134").Value = 14914.0005
$ws.Range("M134").Value = -5602.0905
$ws.Range("N134").Value = -19984.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1132.625
$ws.Range("I31").Value = 1151.4348
$ws.Range("J31").Value = 700
$ws.Range("K31").Value = 1151.4348
$ws.Range("L31").Value = 700
$ws.Range("M31").Value = -856.4348
$ws.Range("N31").Value = -1290

$ws.Range("H34").Value = 1132.625
$ws.Range("I34").Value = 1151.4348
$ws.Range("J34").Value = 700
$ws.Range("K34").Value = 1151.4348
$ws.Range("L34").Value = 700
$ws.Range("M34").Value = -949.4348
$ws.Range("N34").Value = -1104

$ws.Range("H58").Value = 2474.0417
$ws.Range("I58").Value = 1671.2858
$ws.Range("J58").Value = 3597.9
$ws.Range("K58").Value = 1671.2858
$ws.Range("L58").Value = 3597.9
$ws.Range("M58").Value = -1468.2858
$ws.Range("N58").Value = -4003.9

$ws.Range("H132").Value = 5288
$ws.Range("I132").Value = 4888
$ws.Range("J132").Value = 5368
$ws.Range("K132").Value = 14664
$ws.Range("L132").Value = 16104
$ws.Range("M132").Value = -12134
$ws.Range("N132").Value = -21164

$ws.Range("H134").Value = 2581.4062
$ws.Range("I134").Value = 1317.88
$ws.Range("J134").Value = 7094
$ws.Range("K134").Value = 3953.64
$ws.Range("L134").Value = 21282
$ws.Range("M134").Value = -1418.64
$ws.Range("N134").Value = -26352

$ws.Range("H136").Value = 2474.0417
$ws.Range("I136").Value = 1671.2858
$ws.Range("J136").Value = 3597.9
$ws.Range("K136").Value = 5013.857400000001
$ws.Range("L136").Value = 10793.7
$ws.Range("M136").Value = -2463.857400000001
$ws.Range("N136").Value = -15893.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1642.8182
$ws.Range("I5").Value = 920.6667
$ws.Range("J5").Value = 2509.4
$ws.Range("K5").Value = 2762.0001
$ws.Range("L5").Value = 7528.200000000001
$ws.Range("M5").Value = -2650.0001
$ws.Range("N5").Value = -7752.200000000001

$ws.Range("H122").Value = 665.3200000000001
$ws.Range("I122").Value = 279.77777
$ws.Range("J122").Value = 882.1875
$ws.Range("K122").Value = 2517.99993
$ws.Range("L122").Value = 7939.6875
$ws.Range("M122").Value = -67.99992999999995
$ws.Range("N122").Value = -12839.6875

$ws.Range("H135").Value = 1642.8182
$ws.Range("I135").Value = 920.6667
$ws.Range("J135").Value = 2509.4
$ws.Range("K135").Value = 8286.0003
$ws.Range("L135").Value = 22584.6
$ws.Range("M135").Value = -5751.0003
$ws.Range("N135").Value = -27654.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 742
$ws.Range("I107").Value = 1193
$ws.Range("J107").Value = 226.57143
$ws.Range("K107").Value = 1193
$ws.Range("L107").Value = 226.57143
$ws.Range("M107").Value = 727
$ws.Range("N107").Value = -4066.57143

$ws.Range("H132").Value = 2542.6155
$ws.Range("I132").Value = 2010.1875
$ws.Range("J132").Value = 3394.5
$ws.Range("K132").Value = 6030.5625
$ws.Range("L132").Value = 10183.5
$ws.Range("M132").Value = -3500.5625
$ws.Range("N132").Value = -15243.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 734.875
$ws.Range("I93").Value = 705.61536
$ws.Range("J93").Value = 861.6667
$ws.Range("K93").Value = 705.61536
$ws.Range("L93").Value = 861.6667
$ws.Range("M93").Value = 542.38464
$ws.Range("N93").Value = -3357.6667

$ws.Range("H132").Value = 7499.5
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 8332.666999999999
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 24998.001
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -30058.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 50004630
$ws.Range("I132").Value = 83338630
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 250015890
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -250013360
$ws.Range("N132").Value = -15933.5

$ws.Range("H136").Value = 16718451
$ws.Range("I136").Value = 27862478
$ws.Range("J136").Value = 2410
$ws.Range("K136").Value = 83587434
$ws.Range("L136").Value = 7230
$ws.Range("M136").Value = -83584884
$ws.Range("N136").Value = -12330
